# Update column F ("想去人数") values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @(
    @{Row=3;  New=3109},
    @{Row=4;  New=226},
    @{Row=5;  New=119},
    @{Row=6;  New=198},
    @{Row=7;  New=1661},
    @{Row=8;  New=1619},
    @{Row=10; New=366},
    @{Row=12; New=28},
    @{Row=13; New=189},
    @{Row=17; New=230},
    @{Row=21; New=17},
    @{Row=22; New=369},
    @{Row=23; New=179},
    @{Row=24; New=97},
    @{Row=25; New=20},
    @{Row=26; New=19},
    @{Row=27; New=63},
    @{Row=28; New=2086},
    @{Row=29; New=4},
    @{Row=31; New=461},
    @{Row=32; New=197},
    @{Row=36; New=337},
    @{Row=38; New=505},
    @{Row=39; New=9}
)

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates) {
        $ws.Range("F" + $u.Row).Value = $u.New
    }
}
